# Update "想去人数" (F column) counts that changed between data refreshes.
# The same underlying data is duplicated on the "展览" and "全部类型" sheets,
# so both need to be updated identically.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 127
    4  = 1639
    5  = 632
    8  = 11592
    9  = 30
    11 = 459
    12 = 375
    15 = 12415
    16 = 13164
    21 = 248
    22 = 89
    24 = 134
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
